# "Modelo de Difusão Funcionando" - fix the diffusion-model formulas in the
# params sheet: C2 is an independent formula, C3:C6 share formula si="0".
# Both change their divisor from 3 to 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# C2 holds its own (non-shared) formula D2/3 -> D2/5
$ws.Range("C2").Formula = "=D2/5"

# C3:C6 hold the shared formula D#/3 -> D#/5 (anchored at C3, Excel fills
# the relative reference down automatically for C4:C6)
$ws.Range("C3:C6").Formula = "=D3/5"

# Restore the default selection (A1) so the view no longer points at the
# stale E10 cell left over from editing.
$ws.Range("A1").Select()
